$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New word (column B) and frequency (column C) values for rows 2..101,
# reflecting the updated shared-strings table and recomputed counts.
$words = @(
    "셔틀","버스","시간","출발","조","교대","출근","없","아파트","분",
    "노선","리","퇴근","운행","앞","차","정차","사내","사람","차량",
    "문의","회사","관련","정류장","도착","많","인원","이천","푸르","탑승",
    "오","시","이용","서서","경우","부탁","개선","요청","장호원","성남",
    "결제","선","양","홈","송정동","감사","태평","광주","모전","후",
    "건의","카드","타운","수림","쪽","고담","확인","변경","안녕","통상",
    "답변","증차","현재","현금","우","사원","암","늦","일성","곳",
    "좋","추가","안","주공","근무","출퇴근","검토","곤지","정도","생각",
    "방면","오늘","내","방향","동성","통근","마트","증일동","가능","포인트",
    "시내","소요","아침","식권","하이닉스","애","기사님","주차장","야간","힐"
)

$counts = @(
    274,204,151,125,124,111,111,94,93,88,
    69,66,62,62,60,59,57,52,52,50,
    48,47,46,46,46,44,44,44,42,41,
    39,39,38,38,37,35,34,34,33,32,
    30,30,30,30,30,29,29,29,29,28,
    28,28,28,27,27,27,26,26,26,26,
    25,25,24,24,24,23,23,23,23,22,
    22,22,22,22,21,21,21,21,21,20,
    20,20,20,19,19,19,19,19,18,18,
    18,18,17,17,17,17,17,17,17,17
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
    $ws.Cells.Item($row, 3).Value = $counts[$i]
}
